# Weekly price update: insert a new Cilantro price record for
# "Terminal La Palmera de La Serena" at the top of the data block
# (row 111), pushing the existing records down by one row and
# extending the table by one row overall (old last row -> new row 239).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 111; existing rows 111..238 shift to 112..239.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A111").Value = 8
$ws.Range("B111").Value = "Terminal La Palmera de La Serena"
$ws.Range("C111").Value = "Coquimbo"
$ws.Range("D111").Value = 45118
$ws.Range("E111").Value = 4
$ws.Range("F111").Value = 100112040
$ws.Range("G111").Value = "Cilantro"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 2200
$ws.Range("K111").Value = 2500
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 2750
$ws.Range("N111").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O111").Value = "Provincia del Elquí"
$ws.Range("P111").Value = 1833
$ws.Range("Q111").Value = 1.5
$ws.Range("R111").Value = "Hortaliza"
